$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 503. Excel shifts the existing
# rows 503-614 down to 505-616 (standard Rows.Insert behaviour), which is
# exactly the "every later record moved down by two" pattern seen in the
# target diff.
$ws.Rows.Item(503).Insert()
$ws.Rows.Item(503).Insert()

# Populate the two freshly inserted rows (503 and 504) with the new
# weekly price records. Columns follow the same layout as every other
# data row: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion.

$ws.Cells.Item(503, 1).Value = 6
$ws.Cells.Item(503, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(503, 3).Value = "Metropolitana"
$ws.Cells.Item(503, 4).Value = 44951
$ws.Cells.Item(503, 5).Value = 13
$ws.Cells.Item(503, 6).Value = 100112052
$ws.Cells.Item(503, 7).Value = "Albahaca"
$ws.Cells.Item(503, 8).Value = "Sin especificar"
$ws.Cells.Item(503, 9).Value = "Primera"
$ws.Cells.Item(503, 10).Value = 1080
$ws.Cells.Item(503, 11).Value = 2500
$ws.Cells.Item(503, 12).Value = 3000
$ws.Cells.Item(503, 13).Value = 2708
$ws.Cells.Item(503, 14).Value = "`$/docena de matas"
$ws.Cells.Item(503, 15).Value = "Región Metropolitana"
$ws.Cells.Item(503, 16).Value = 451
$ws.Cells.Item(503, 17).Value = 6
$ws.Cells.Item(503, 18).Value = "Hortaliza"

$ws.Cells.Item(504, 1).Value = 6
$ws.Cells.Item(504, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(504, 3).Value = "Metropolitana"
$ws.Cells.Item(504, 4).Value = 44951
$ws.Cells.Item(504, 5).Value = 13
$ws.Cells.Item(504, 6).Value = 100112052
$ws.Cells.Item(504, 7).Value = "Albahaca"
$ws.Cells.Item(504, 8).Value = "Sin especificar"
$ws.Cells.Item(504, 9).Value = "Segunda"
$ws.Cells.Item(504, 10).Value = 450
$ws.Cells.Item(504, 11).Value = 2000
$ws.Cells.Item(504, 12).Value = 2000
$ws.Cells.Item(504, 13).Value = 2000
$ws.Cells.Item(504, 14).Value = "`$/docena de matas"
$ws.Cells.Item(504, 15).Value = "Región Metropolitana"
$ws.Cells.Item(504, 16).Value = 333
$ws.Cells.Item(504, 17).Value = 6
$ws.Cells.Item(504, 18).Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of the
# column (YYYY-MM-DD HH:MM:SS), matching style index 2 used throughout.
$ws.Range("D503:D504").NumberFormat = "YYYY-MM-DD HH:MM:SS"
